# Apply the daily cryptos-list refresh (prices + 1h volume deltas) scraped
# by the GitHub Actions job. A handful of rows also got re-ranked, so their
# Coin/Link/Price/Volume columns move to a different row.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Plain strings can be written directly; values that *look* like a bare
# number (e.g. "553.66") need a leading apostrophe so Excel keeps them as
# text instead of silently converting the cell to a Number, and the style
# is reset afterwards so the quote-prefix flag does not linger as a style
# change on the cell.
function Set-CellText($cellRef, $value) {
    $cell = $ws.Range($cellRef)
    $cell.Value = "'" + $value
    $cell.Style = "Normal"
}

$updates = @(
    ,@("D2", "62.941.55")
    ,@("E2", "  -1.82%  ")
    ,@("D3", "2.679.92")
    ,@("E3", "  -2.29%  ")
    ,@("E4", "  -0.03%  ")
    ,@("D5", "553.66")
    ,@("E5", "  -3.20%  ")
    ,@("D6", "158.28")
    ,@("E6", "  -1.61%  ")
    ,@("E7", "  +0.02%  ")
    ,@("E8", "  -0.64%  ")
    ,@("E9", "  -3.67%  ")
    ,@("E10", "  -1.81%  ")
    ,@("D11", "0.367")
    ,@("E11", "  -4.53%  ")
    ,@("D12", "5.38")
    ,@("E12", "  -7.21%  ")
    ,@("D13", "3.152.45")
    ,@("E13", "  -2.31%  ")
    ,@("D14", "26.23")
    ,@("E14", "  -2.50%  ")
    ,@("D15", "62.788.58")
    ,@("E15", "  -1.79%  ")
    ,@("E16", "  -2.99%  ")
    ,@("D17", "2.680.87")
    ,@("E17", "  -2.43%  ")
    ,@("D18", "11.87")
    ,@("E18", "  -2.60%  ")
    ,@("D19", "4.61")
    ,@("E19", "  -4.34%  ")
    ,@("D20", "344.79")
    ,@("E20", "  -2.88%  ")
    ,@("E21", "  -4.91%  ")
    ,@("D22", "0.999")
    ,@("E22", "  -0.04%  ")
    ,@("E23", "  -3.19%  ")
    ,@("D24", "63.37")
    ,@("E24", "  -1.55%  ")
    ,@("E25", "  -1.94%  ")
    ,@("E26", "  +0.02%  ")
    ,@("D27", "8.19")
    ,@("E27", "  -3.46%  ")
    ,@("D28", "0.0₃0857")
    ,@("E28", "  -6.88%  ")
    ,@("D29", "1.38")
    ,@("E29", "  +2.09%  ")
    ,@("D30", "7.23")
    ,@("E30", "  -0.47%  ")
    ,@("E31", "  -1.73%  ")
    ,@("D32", "165.06")
    ,@("E32", "  +0.48%  ")
    ,@("E33", "  -2.04%  ")
    ,@("E34", "  -0.33%  ")
    ,@("E35", "  +0.02%  ")
    ,@("D36", "19.50")
    ,@("E36", "  -3.14%  ")
    ,@("D37", "1.78")
    ,@("E37", "  -2.12%  ")
    ,@("D38", "348.89")
    ,@("E38", "  -0.28%  ")
    ,@("B39", "SuiNetwork")
    ,@("C39", "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui")
    ,@("D39", "0.961")
    ,@("E39", "  -3.60%  ")
    ,@("B40", "RenderToken")
    ,@("C40", "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr")
    ,@("D40", "6.32")
    ,@("E40", "  -1.26%  ")
    ,@("D41", "4.00")
    ,@("E41", "  -3.03%  ")
    ,@("D42", "38.30")
    ,@("E42", "  -1.01%  ")
    ,@("B43", "EnergySwap")
    ,@("C43", "https://coinranking.com/coin/SbWqqTui-+energyswap-ens")
    ,@("D43", "20.37")
    ,@("E43", "  -4.22%  ")
    ,@("B44", "InjectiveProtocol")
    ,@("C44", "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj")
    ,@("D44", "20.74")
    ,@("E44", "  -6.11%  ")
    ,@("B45", "Hedera")
    ,@("C45", "https://coinranking.com/coin/jad286TjB+hedera-hbar")
    ,@("D45", "0.0561")
    ,@("E45", "  -4.19%  ")
    ,@("B46", "Mantle")
    ,@("C46", "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt")
    ,@("D46", "0.616")
    ,@("E46", "  -1.81%  ")
    ,@("D47", "0.998")
    ,@("E47", "  -0.02%  ")
    ,@("E48", "  -0.04%  ")
    ,@("D49", "0.0972")
    ,@("E49", "  -3.56%  ")
    ,@("E50", "  -3.71%  ")
    ,@("D51", "2.098.77")
    ,@("E51", "  -1.98%  ")
)

foreach ($u in $updates) {
    Set-CellText $u[0] $u[1]
}
